# Update the "Metadata" sheet (sheet1): version, date, publisher, new
# Jurisdiction row, and the resulting shift of all subsequent rows caused
# by removing the old duplicated "Contact" row.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

$meta.Range("B9").Value = "Alvearie Team"

$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

$meta.Range("A11").Value = "Description"
$meta.Range("B11").Value = "ID or url reference to template to be used for communication"

$meta.Range("A12").Value = "Purpose"
$meta.Range("B12").Value = $null

$meta.Range("A13").Value = "Copyright"
$meta.Range("B13").Value = $null

$meta.Range("A14").Value = "FHIR Version"
$meta.Range("B14").Value = "4.0.1"

$meta.Range("A15").Value = "Kind"
$meta.Range("B15").Value = "complex-type"

$meta.Range("A16").Value = "Type"
$meta.Range("B16").Value = "Extension"

$meta.Range("A17").Value = "Base Definition"
$meta.Range("B17").Value = "http://hl7.org/fhir/StructureDefinition/Extension"

$meta.Range("A18").Value = "Abstract"
$meta.Range("B18").Value = "'false"

$meta.Range("A19").Value = "Derivation"
$meta.Range("B19").Value = "constraint"

$meta.Range("A20").Value = "Context"
$meta.Range("B20").Value = "element:Element"

# The old table had 21 rows; the new one only has 20 (the duplicated
# "Contact" row is gone), so drop the now-superfluous last row.
$meta.Rows.Item(21).Delete()

# Update the "Elements" sheet (sheet2): the root Extension row's Short /
# Definition now mirror the template's Title / Description instead of the
# generic "Extension" / "An Extension" text.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Communication Template"
$elements.Range("L2").Value = "ID or url reference to template to be used for communication"
